$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.796.69"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.635.69"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  -0.08%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "215.42"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.08%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "19.86"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +1.56%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0787"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("D12").Value = "1.654.94"
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Value = "1.860.54"
$ws.Range("E14").Value = "  -0.08%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.556"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").Value = "0.0₃0776"
$ws.Range("E16").Value = "  +2.28%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "63.10"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "25.813.23"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("E20").Value = "  +2.82%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "194.45"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("E22").Value = "  +0.65%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "6.16"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +1.00%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("E25").Value = "  -1.76%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "139.80"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -5.12%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "6.83"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("E31").Value = "  +1.57%  "
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("E33").Value = "  +1.08%  "
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("E35").Value = "  +0.38%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.898"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -0.48%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").Value = "1.113.03"
$ws.Range("E39").Value = "  -1.60%  "
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("E42").Value = "  +1.05%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "99.29"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +1.56%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.800"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("E45").Value = "  -0.91%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "55.59"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +0.48%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.51"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +13.13%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "7.73"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("E49").Value = "  -4.92%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.0502"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("E51").Value = "  -0.37%  "
